$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.775.24"
$ws.Range("E2").Value = "  +1.97%  "

$ws.Range("D3").Value = "1.876.66"
$ws.Range("E3").Value = "  +2.16%  "

$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").Value = "'326.77"
$ws.Range("E5").Value = "  -1.00%  "

$ws.Range("D6").Value = "'1.003"
$ws.Range("E6").Value = "  +0.06%  "

$ws.Range("D7").Value = "'0.4650"
$ws.Range("E7").Value = "  +1.13%  "

$ws.Range("D8").Value = "'0.3923"
$ws.Range("E8").Value = "  +1.75%  "

$ws.Range("D9").Value = "'0.07925"
$ws.Range("E9").Value = "  +0.87%  "

$ws.Range("D10").Value = "'0.9750"
$ws.Range("E10").Value = "  +1.71%  "

$ws.Range("D11").Value = "'22.34"
$ws.Range("E11").Value = "  +1.96%  "

$ws.Range("D12").Value = "1.840.19"
$ws.Range("E12").Value = "  +0.19%  "

$ws.Range("D13").Value = "'5.754"
$ws.Range("E13").Value = "  +0.78%  "

$ws.Range("D14").Value = "'6.963"
$ws.Range("E14").Value = "  +0.80%  "

$ws.Range("D15").Value = "'0.07025"
$ws.Range("E15").Value = "  +2.74%  "

$ws.Range("D16").Value = "'88.60"
$ws.Range("E16").Value = "  +1.98%  "

$ws.Range("D17").Value = "'1.004"
$ws.Range("E17").Value = "  +0.12%  "

$ws.Range("D18").Value = "'0.00001010"
$ws.Range("E18").Value = "  +1.66%  "

$ws.Range("D19").Value = "'16.98"
$ws.Range("E19").Value = "  +0.37%  "

$ws.Range("D21").Value = "28.775.41"
$ws.Range("E21").Value = "  +1.88%  "

$ws.Range("D22").Value = "'5.340"
$ws.Range("E22").Value = "  +0.10%  "

$ws.Range("D23").Value = "'11.12"
$ws.Range("E23").Value = "  +1.37%  "

$ws.Range("D25").Value = "2.073.34"
$ws.Range("E25").Value = "  +0.97%  "

$ws.Range("D26").Value = "'153.54"
$ws.Range("E26").Value = "  +0.32%  "

$ws.Range("D27").Value = "'19.41"
$ws.Range("E27").Value = "  +1.08%  "

$ws.Range("D28").Value = "'5.754"
$ws.Range("E28").Value = "  +0.93%  "

$ws.Range("D29").Value = "'2.009"
$ws.Range("E29").Value = "  +1.60%  "

$ws.Range("D30").Value = "'119.74"
$ws.Range("E30").Value = "  +2.52%  "

$ws.Range("D31").Value = "'0.09380"
$ws.Range("E31").Value = "  +1.04%  "

$ws.Range("D32").Value = "'0.9396"
$ws.Range("E32").Value = "  -0.16%  "

$ws.Range("D33").Value = "'5.335"
$ws.Range("E33").Value = "  +1.14%  "

$ws.Range("D34").Value = "'1.351"
$ws.Range("E34").Value = "  +2.12%  "

$ws.Range("D35").Value = "'3.354"
$ws.Range("E35").Value = "  -2.70%  "

$ws.Range("D36").Value = "'0.05869"
$ws.Range("E36").Value = "  -2.15%  "

$ws.Range("D37").Value = "'0.02127"
$ws.Range("E37").Value = "  -1.04%  "

$ws.Range("D39").Value = "'7.941"
$ws.Range("E39").Value = "  +4.29%  "

$ws.Range("D40").Value = "'0.5674"
$ws.Range("E40").Value = "  +1.04%  "

$ws.Range("D43").Value = "'0.07245"
$ws.Range("E43").Value = "  +3.25%  "

$ws.Range("D44").Value = "'11.72"
$ws.Range("E44").Value = "  +0.64%  "

$ws.Range("D45").Value = "'0.5328"
$ws.Range("E45").Value = "  +0.96%  "

$ws.Range("D46").Value = "'2.143"
$ws.Range("E46").Value = "  -4.64%  "

$ws.Range("D47").Value = "'1.136"
$ws.Range("E47").Value = "  -7.10%  "

$ws.Range("D48").Value = "'1.855"
$ws.Range("E48").Value = "  +1.30%  "

$ws.Range("D49").Value = "'113.78"
$ws.Range("E49").Value = "  +0.90%  "

$ws.Range("D50").Value = "'2.352"
$ws.Range("E50").Value = "  +0.91%  "

$ws.Range("D51").Value = "'1.003"
$ws.Range("E51").Value = "  +0.10%  "

$ws.Range("E20").Value = "  +0.41%  "

$ws.Range("E24").Value = "  -0.87%  "

$ws.Range("E38").Value = "  +0.13%  "

# Row 41 and 42: Aptos/Algorand swap position with updated values
$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").Value = "'0.1791"
$ws.Range("E41").Value = "  +1.01%  "

$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").Value = "'9.968"
$ws.Range("E42").Value = "  -0.10%  "
